# Add self assessment data to row 10 (Jonas Antunes), both the left block
# (C:G) and the mirrored right block (L:P).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10 previously held the placeholder text "118----"; it becomes a formula
# that evaluates to the numeric id, same as the other rows in the table.
$ws.Range("C10").Formula = "=1181478"
$ws.Range("L10").Formula = "=1181478"

# F10/G10 previously were blank; fill in the self assessment ratings.
$ws.Range("F10").Value = "good"
$ws.Range("G10").Value = "very good"
